$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds free-form text values straight from the
# source feed (e.g. "62.715.53", "1.00", "0.0000100"); some of those look
# like numbers to the input parser and would otherwise get silently
# coerced into numeric cells (dropping trailing zeros, etc). Mark the whole
# column range as Text up front so every write below lands as a string,
# then restore the default Normal style once all writes are done so no
# stray formatting is left on the cells.
$priceCol = $ws.Range("D2:D51")
$priceCol.NumberFormat = "@"

$ws.Range("D2").Value = '62.715.53'
$ws.Range("E2").Value = '  +0.82%  '
$ws.Range("D3").Value = '2.443.03'
$ws.Range("E3").Value = '  +0.98%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '567.38'
$ws.Range("E5").Value = '  +0.56%  '
$ws.Range("D6").Value = '146.12'
$ws.Range("E6").Value = '  +2.15%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("E8").Value = '  +1.00%  '
$ws.Range("E9").Value = '  +2.12%  '
$ws.Range("D11").Value = '5.27'
$ws.Range("E11").Value = '  -1.27%  '
$ws.Range("E12").Value = '  +0.00%  '
$ws.Range("D13").Value = '0.0000186'
$ws.Range("E13").Value = '  +5.62%  '
$ws.Range("D14").Value = '26.93'
$ws.Range("E14").Value = '  +4.63%  '
$ws.Range("D15").Value = '2.834.02'
$ws.Range("E15").Value = '  -0.74%  '
$ws.Range("D16").Value = '62.537.24'
$ws.Range("E16").Value = '  +0.79%  '
$ws.Range("D17").Value = '2.432.40'
$ws.Range("E17").Value = '  +0.76%  '
$ws.Range("D18").Value = '11.30'
$ws.Range("E18").Value = '  -0.44%  '
$ws.Range("D19").Value = '6.95'
$ws.Range("E19").Value = '  +1.24%  '
$ws.Range("D20").Value = '325.24'
$ws.Range("E20").Value = '  +0.46%  '
$ws.Range("E21").Value = '  +0.06%  '
$ws.Range("D22").Value = '1.00'
$ws.Range("E22").Value = '  -0.07%  '
$ws.Range("D23").Value = '67.44'
$ws.Range("E23").Value = '  +2.20%  '
$ws.Range("D24").Value = '1.75'
$ws.Range("E24").Value = '  +2.43%  '
$ws.Range("D25").Value = '8.76'
$ws.Range("E25").Value = '  -2.19%  '
$ws.Range("D26").Value = '0.0000100'
$ws.Range("E26").Value = '  +5.21%  '
$ws.Range("D27").Value = '561.07'
$ws.Range("E27").Value = '  -3.50%  '
$ws.Range("D28").Value = '2.562.45'
$ws.Range("E28").Value = '  +1.24%  '
$ws.Range("E29").Value = '  -0.23%  '
$ws.Range("D30").Value = '8.37'
$ws.Range("E30").Value = '  +1.50%  '
$ws.Range("E31").Value = '  +1.55%  '
$ws.Range("E32").Value = '  -0.61%  '
$ws.Range("E33").Value = '  +0.30%  '
$ws.Range("E34").Value = '  +0.60%  '
$ws.Range("D35").Value = '4.90'
$ws.Range("E35").Value = '  +3.61%  '
$ws.Range("D36").Value = '0.998'
$ws.Range("E36").Value = '  -0.17%  '
$ws.Range("D37").Value = '0.384'
$ws.Range("E37").Value = '  +0.46%  '
$ws.Range("D38").Value = '5.46'
$ws.Range("E38").Value = '  -1.69%  '
$ws.Range("D39").Value = '18.86'
$ws.Range("E39").Value = '  +0.84%  '
$ws.Range("D40").Value = '150.24'
$ws.Range("E40").Value = '  -1.15%  '
$ws.Range("D41").Value = '1.83'
$ws.Range("E41").Value = '  +1.41%  '
$ws.Range("E42").Value = '  +0.56%  '
$ws.Range("D43").Value = '2.42'
$ws.Range("D44").Value = '149.13'
$ws.Range("E44").Value = '  +0.25%  '
$ws.Range("E45").Value = '  +1.00%  '
$ws.Range("E46").Value = '  +0.48%  '
$ws.Range("D47").Value = '20.59'
$ws.Range("E47").Value = '  +2.51%  '
$ws.Range("E48").Value = '  +1.22%  '
$ws.Range("E49").Value = '  +1.09%  '
$ws.Range("E50").Value = '  +1.88%  '
$ws.Range("D51").Value = '11.60'
$ws.Range("E51").Value = '  +0.45%  '

# Restore normal formatting on the price column now that the text values are set.
$priceCol.Style = "Normal"
